$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) whose new value is numeric-looking. Excel would silently
# coerce these into floating point numbers (losing exact formatting / introducing
# binary rounding noise), so mark them as Text before writing the literal string.
$forceTextCells = @('D4', 'D5', 'D6', 'D13', 'D14', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D27', 'D28', 'D33', 'D37', 'D43', 'D44', 'D45', 'D47', 'D49', 'D50', 'D51')
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.752.73'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '3.848.11'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '602.03'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').Value = '170.83'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('D7').Value = '3.846.64'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').Value = '0.0000285'
$ws.Range('E13').Value = '  +14.40%  '
$ws.Range('D14').Value = '37.11'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').Value = '4.498.04'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '3.796.82'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').Value = '68.768.11'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '18.32'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = '7.39'
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('D20').Value = '0.112'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').Value = '11.18'
$ws.Range('E21').Value = '  +4.09%  '
$ws.Range('D22').Value = '473.06'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('D24').Value = '0.0000165'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Value = '12.14'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = '10.37'
$ws.Range('E28').Value = '  +3.32%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').Value = '4.000.22'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('D33').Value = '31.41'
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('E34').Value = '  -0.46%  '
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('D36').Value = '3.815.97'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '3.90'
$ws.Range('E37').Value = '  +18.12%  '
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').Value = '0.317'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('B44').Value = 'FLOKI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D44').Value = '0.000303'
$ws.Range('E44').Value = '  +10.21%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '2.00'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '420.70'
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('D49').Value = '46.37'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('D50').Value = '141.58'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('D51').Value = '0.0358'
$ws.Range('E51').Value = '  -0.28%  '

# Reset number format on the forced-text cells back to the default (General/Normal)
# style so no stray cell formatting is left behind - only the literal text remains.
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).Style = "Normal"
}
